# Applies the text corrections described in the commit diff to
# worksheet 'Tab_4a_Indikatorenblätter' of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab_4a_Indikatorenblätter")

$text = @'
Ausgehend von der Zielformulierung würde bei Fortsetzung der Entwicklung der letzten sechs Jahre der gesetzte Zielwert bereits deutlich früher (im Berichtsjahr 2024) unterschritten, sodass der Indikator 2.1.a für das Jahr 2022 mit „Sonne“ bewertet wird.
<u>Hinweis:</u> Der Indikator wird als gleitender Fünfjahresdurchschnitt dargestellt, . h. der für das Zieljahr 2030 relevante Wert wird aus den Einzelwerten der Jahre 2026 bis 2030 berechnet.
'@
$ws.Range("L3").Value2 = $text

$text = @'
Es handelt sich um einen neuen Indikator aus der Weiterentwicklung 2025 der Deutschen Nachhaltigkeitsstrategie.
'@
$ws.Range("F12").Value2 = $text

$text = @'
Es handelt sich um einen neuen Indikator aus der Weiterentwicklung 2025 der Deutschen Nachhaltigkeitsstrategie.
'@
$ws.Range("F15").Value2 = $text

$text = @'
Ausgehend von der Zielformulierung deutet beim Indikator 4.2.a die durchschnittliche Entwicklung der letzten sechs Jahre (trotz Stagnation in den letzten Jahren) in die richtige Richtung. Bei Fortsetzung dieser Entwicklung wird der Indikator weiterhin geringfügig ansteigen, das gesetzte Ziel in 2030 aber weit verfehlen, sodass der Indikator 4.2.a für das Jahr 2024 mit „Wolke“ bewertet wird. Dagegen hat sich der Anteil beim Indikator 4.2.b zwischen 2019 und 2024 sogar leicht verringert und damit nicht in die gewünschte Richtung entwickelt. Der Indikator 4.2.b wird für das Jahr 2024 mit „Gewitter“ bewertet.
'@
$ws.Range("L16").Value2 = $text

$text = @'
Es handelt sich um einen neuen Indikator aus der Weiterentwicklung 2025 der Deutschen Nachhaltigkeitsstrategie.
'@
$ws.Range("F17").Value2 = $text

$text = @'
Es handelt sich um einen neuen Indikator aus der Weiterentwicklung 2025 der Deutschen Nachhaltigkeitsstrategie.
'@
$ws.Range("F18").Value2 = $text

$text = @'
Es handelt sich um einen neuen Indikator aus der Weiterentwicklung 2025 der Deutschen Nachhaltigkeitsstrategie.
'@
$ws.Range("F22").Value2 = $text

$text = @'
Es handelt sich um einen neuen Indikator aus der Weiterentwicklung 2025 der Deutschen Nachhaltigkeitsstrategie.
'@
$ws.Range("F24").Value2 = $text

$text = @'
Ausgehend von der Zielformulierung ist für die Bewertung des Indikators nicht relevant, ob die Steigerung durch eine Steigerung des Zählers oder eine Verringerung des Nenners erreicht wird. Die Werte des Indikators 8.3 sind sowohl im Jahr 2024, als auch im Durchschnitt der Jahre 2019 bis 2024 gesunken, d. h. die Werte entwickelten sich nicht in die gewünschte Richtung. Der Indikator 8.3 wird für das Jahr 2024 mit „Gewitter“ bewertet.
'@
$ws.Range("L34").Value2 = $text

$text = @'
Es handelt sich um einen neuen Indikator aus der Weiterentwicklung 2025 der Deutschen Nachhaltigkeitsstrategie.
'@
$ws.Range("F37").Value2 = $text

$text = @'
Es handelt sich um einen neuen Indikator aus der Weiterentwicklung 2025 der Deutschen Nachhaltigkeitsstrategie.
'@
$ws.Range("F38").Value2 = $text

$text = @'
Ausgehend von der Zielformulierung wird für jedes Jahr die Differenz zwischen dem EU-Wert und dem Wert für Deutschland gebildet. Für den Indikator 10.2 werden (aufgrund methodischer Änderungen am Erhebungskonzept) die Indikatorwerte von 2020 bis 2023 betrachtet. Die Differenz für das Jahr 2023 ist positiv, d. h. der Koeffizient in Deutschland ist niedriger als der EU-Koeffizient. Das Ziel für diesen Indikator ist erfüllt. Da sich auch die Differenz in den letzten vier Jahren im Durchschnitt vergrößert hat, wird der Indikator 10.2 für das Jahr 2023 mit „Sonne“ bewertet.
<u>Hinweis:</u> Der Bezug auf den EU-Koeffizienten als Zielgröße bedeutet, dass Indikatoren positiv bewertet werden können, auch wenn sich der Gini-Koeffizient in Deutschland negativ entwickelt. Außerdem ist bei dem Indikator festzustellen, dass durch die Verläufe der beiden Koeffizienten auf ähnlich hohem Niveau und ohne eine deutliche steigende oder sinkende Tendenz, sowohl die Differenz zwischen deutschem und europäischem Wert als auch die Richtung der durchschnittlichen Entwicklung des deutschen Koeffizienten starken Schwankungen unterworfen sind, sodass die Bewertungen schon durch geringfügige Änderungen in dem Indikator stark beeinflusst wird.
'@
$ws.Range("L42").Value2 = $text

$text = @'
Für den Indikator 11.3.a werden (aufgrund methodischer Änderungen am Erhebungskonzept) die Indikatorwerte von 2020 bis 2023 betrachtet. In den vergangenen vier Jahren lag der Wert des Indikators unter dem Zielwert von 13 %. Das Ziel war erfüllt. Aufgrund der zugrundeliegenden Bewertungsmethodik muss das Ziel nun in jedem Folgejahr gehalten werden und die durchschnittliche Veränderung darf nicht in Richtung einer Verschlechterung weisen. Im Jahr 2023 lag der Wert des Indikators genau auf dem Zielwert von 13 %. Dieses Teilziel gilt als erfüllt. Da aber die durchschnittliche Entwicklung der letzten vier Jahre nicht in die gewünschte Richtung weist, wird der Indikator 11.3.a für das Jahr 2023 mit "Leicht bewölkt" bewertet.
'@
$ws.Range("L48").Value2 = $text

$text = @'
Es handelt sich um einen neuen Indikator aus der Weiterentwicklung 2025 der Deutschen Nachhaltigkeitsstrategie.
'@
$ws.Range("F49").Value2 = $text

$text = @'
Es handelt sich um einen neuen Indikator aus der Weiterentwicklung 2025 der Deutschen Nachhaltigkeitsstrategie.
'@
$ws.Range("F53").Value2 = $text

$text = @'
12.3.a: Keine Bewertung möglich. Zu wenig Datenpunkte.
12.3.b: Ausgehend von der Zielformulierung entwickelte sich Indikator 12.3.b im letzten Jahr (2022) in die gewünschte Richtung. Die Entwicklung im Durchschnitt der letzten sechs Jahre ging jedoch nicht in die gewünschte Richtung, sodass der Indikator 12.3.b für das Jahr 2022 mit "Wolke" bewertet wird.
12.3.c: Keine Bewertung möglich. Zu wenig Datenpunkte.
'@
$ws.Range("L54").Value2 = $text

$text = @'
Es handelt sich um einen neuen Indikator aus der Weiterentwicklung 2025 der Deutschen Nachhaltigkeitsstrategie.
'@
$ws.Range("F55").Value2 = $text

$text = @'
Für Indikator 13.1.a sind verschiedene Ziele definiert, die in unterschiedlichen Jahren erreicht werden sollen. Für die Bewertung des Indikators ist das am nächsten in der Zukunft liegende Ziel relevant. Bei Beibehaltung der Entwicklung der letzten sechs Jahre wird der gesetzte Zielwert in 2030 erreicht, sodass der Indikator 13.1.a für das Jahr 2023 mit "Sonne" bewertet wird.
'@
$ws.Range("L56").Value2 = $text

$text = @'
Es handelt sich um einen neuen Indikator aus der Weiterentwicklung 2025 der Deutschen Nachhaltigkeitsstrategie.
'@
$ws.Range("F62").Value2 = $text

$text = @'
Es handelt sich um einen neuen Indikator aus der Weiterentwicklung 2025 der Deutschen Nachhaltigkeitsstrategie.
'@
$ws.Range("F63").Value2 = $text

$text = @'
Es handelt sich um einen neuen Indikator aus der Weiterentwicklung 2025 der Deutschen Nachhaltigkeitsstrategie.
'@
$ws.Range("F64").Value2 = $text

$text = @'
Es handelt sich um einen neuen Indikator aus der Weiterentwicklung 2025 der Deutschen Nachhaltigkeitsstrategie.
'@
$ws.Range("F66").Value2 = $text
